$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-08-13 Wednesday"; new="2025-08-14 Thursday"},
    @{old="23×51=1173"; new="12×29=348"},
    @{old="11×93=1023"; new="55×64=3520"},
    @{old="25×67=1675"; new="52×69=3588"},
    @{old="53×86=4558"; new="27×29=783"},
    @{old="15×32=480"; new="84×79=6636"},
    @{old="13×26=338"; new="23×93=2139"},
    @{old="72×46=3312"; new="71×42=2982"},
    @{old="71×73=5183"; new="47×33=1551"},
    @{old="82×91=7462"; new="84×87=7308"},
    @{old="76×50=3800"; new="26×89=2314"},
    @{old="46×61=2806"; new="13×60=780"},
    @{old="19×80=1520"; new="28×44=1232"},
    @{old="34×38=1292"; new="47×64=3008"},
    @{old="34×32=1088"; new="42×64=2688"},
    @{old="45×21=945"; new="97×99=9603"},
    @{old="12×97=1164"; new="47×87=4089"},
    @{old="83×67=5561"; new="85×65=5525"},
    @{old="97×50=4850"; new="39×17=663"},
    @{old="95×76=7220"; new="99×24=2376"},
    @{old="42×55=2310"; new="50×53=2650"},
    @{old="15×21=315"; new="39×42=1638"},
    @{old="54×58=3132"; new="45×13=585"},
    @{old="94×80=7520"; new="41×85=3485"},
    @{old="81×13=1053"; new="71×57=4047"},
    @{old="55×98=5390"; new="69×80=5520"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2) | Out-Null
}
